$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns (D, E) hold numeric- and percent-looking text values that
# must stay as text (matching the original inlineStr cells), so force a Text number
# format on them before writing so Excel does not auto-convert them to numbers.
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "D7", "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12",
    "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17",
    "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22",
    "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D38", "E38",
    "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43",
    "D44", "E44", "D46", "E46", "E47", "D48", "D49", "E49", "D50", "E50",
    "D51", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values.
$ws.Range("D2").Value = "331.56"
$ws.Range("E2").Value = "0.45%"
$ws.Range("D3").Value = "45.46"
$ws.Range("E3").Value = "2.99%"
$ws.Range("D4").Value = "5.619"
$ws.Range("E4").Value = "2.31%"
$ws.Range("D5").Value = "0.08341"
$ws.Range("E5").Value = "4.21%"
$ws.Range("D6").Value = "2.096"
$ws.Range("E6").Value = "6.14%"
$ws.Range("D7").Value = "0.9627"
$ws.Range("E7").Value = "1.19%"
$ws.Range("E8").Value = "-0.92%"
$ws.Range("D9").Value = "0.1163"
$ws.Range("E9").Value = "5.81%"
$ws.Range("D10").Value = "0.1926"
$ws.Range("E10").Value = "0.71%"
$ws.Range("D11").Value = "10.36"
$ws.Range("E11").Value = "-0.81%"
$ws.Range("D12").Value = "0.09930"
$ws.Range("E12").Value = "-0.18%"
$ws.Range("D13").Value = "0.04612"
$ws.Range("E13").Value = "-3.70%"
$ws.Range("D14").Value = "0.1061"
$ws.Range("E14").Value = "-0.33%"
$ws.Range("D15").Value = "0.001291"
$ws.Range("E15").Value = "1.33%"
$ws.Range("D16").Value = "0.006096"
$ws.Range("E16").Value = "2.54%"
$ws.Range("E17").Value = "0.19%"
$ws.Range("D18").Value = "4.438"
$ws.Range("E18").Value = "1.07%"
$ws.Range("D19").Value = "0.3342"
$ws.Range("E19").Value = "-4.26%"
$ws.Range("D20").Value = "0.1393"
$ws.Range("E20").Value = "-1.90%"
$ws.Range("D21").Value = "0.2880"
$ws.Range("E21").Value = "11.33%"
$ws.Range("D22").Value = "0.04179"
$ws.Range("E22").Value = "2.27%"
$ws.Range("D23").Value = "0.001317"
$ws.Range("E23").Value = "3.47%"
$ws.Range("D24").Value = "0.004557"
$ws.Range("E24").Value = "4.18%"
$ws.Range("D25").Value = "0.0001304"
$ws.Range("E25").Value = "8.61%"
$ws.Range("D26").Value = "0.0003750"
$ws.Range("E26").Value = "0.14%"
$ws.Range("D38").Value = "0.02713"
$ws.Range("E38").Value = "4.92%"
$ws.Range("D39").Value = "0.05762"
$ws.Range("E39").Value = "1.20%"
$ws.Range("D40").Value = "0.007884"
$ws.Range("E40").Value = "4.25%"
$ws.Range("D41").Value = "0.1435"
$ws.Range("E41").Value = "2.51%"
$ws.Range("D42").Value = "0.007268"
$ws.Range("E42").Value = "-1.23%"
$ws.Range("D43").Value = "0.002026"
$ws.Range("E43").Value = "0.50%"
$ws.Range("D44").Value = "0.009119"
$ws.Range("E44").Value = "9.15%"
$ws.Range("D46").Value = "0.00007113"
$ws.Range("E46").Value = "-0.27%"
$ws.Range("E47").Value = "0.25%"
$ws.Range("D48").Value = "0.0005816"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "0.003508"
$ws.Range("E49").Value = "-0.66%"
$ws.Range("B50").Value = "BOLO"
$ws.Range("C50").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D50").Value = "0.003527"
$ws.Range("E50").Value = "-0.82%"
$ws.Range("D51").Value = "0.00002106"
$ws.Range("E51").Value = "0.25%"
